# Restore the deck's theme color palette back to the stock "Office Theme"
# scheme (the file had drifted to the "Integral" theme's greens/yellows).
# Font scheme / format scheme are already identical between the two themes,
# so only the 12 theme colors (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink) need
# to be put back to their original "Office" values.

function HexToRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    # PowerPoint's ColorFormat.RGB uses 0x00BBGGRR ordering.
    return $r -bor ($g -shl 8) -bor ($b -shl 16)
}

$officeTheme = @{
    1  = "000000"  # dk1
    2  = "FFFFFF"  # lt1
    3  = "44546A"  # dk2
    4  = "E7E6E6"  # lt2
    5  = "5B9BD5"  # accent1
    6  = "ED7D31"  # accent2
    7  = "A5A5A5"  # accent3
    8  = "FFC000"  # accent4
    9  = "4472C4"  # accent5
    10 = "70AD47"  # accent6
    11 = "0563C1"  # hlink
    12 = "954F72"  # folHlink
}

$p = $ppt.ActivePresentation

# The theme's ColorScheme is shared by every slide/master in the deck, so
# touching it through any one slide updates the underlying theme part.
$cs = $p.Slides.Item(1).ColorScheme

foreach ($idx in $officeTheme.Keys) {
    $cs.Colors($idx).RGB = HexToRgb $officeTheme[$idx]
}
